$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.1820580474934037
$ws.Range("C2").Value = 0.5593667546174143
$ws.Range("J2").Value = 0.0158311345646438
$ws.Range("P2").Value = 0.1424802110817942
$ws.Range("S2").Value = 0.1002638522427441
$ws.Range("B3").Value = 0.008658008658008658
$ws.Range("C3").Value = 0.01731601731601732
$ws.Range("J3").Value = 0.04329004329004329
$ws.Range("P3").Value = 0.7575757575757576
$ws.Range("S3").Value = 0.1731601731601732
$ws.Range("J4").Value = 0.119047619047619
$ws.Range("P4").Value = 0.6428571428571429
$ws.Range("S4").Value = 0.2380952380952381
$ws.Range("B6").Value = 0.07446808510638298
$ws.Range("D6").Value = 0.003546099290780142
$ws.Range("F6").Value = 0.07446808510638298
$ws.Range("J6").Value = 0.2588652482269503
$ws.Range("O6").Value = 0.02836879432624113
$ws.Range("Q6").Value = 0.1170212765957447
$ws.Range("R6").Value = 0.0851063829787234
$ws.Range("S6").Value = 0.3581560283687943
$ws.Range("B7").Value = 0.1681034482758621
$ws.Range("D7").Value = 0.02586206896551724
$ws.Range("F7").Value = 0.02155172413793104
$ws.Range("J7").Value = 0.1379310344827586
$ws.Range("O7").Value = 0.03448275862068965
$ws.Range("Q7").Value = 0.1810344827586207
$ws.Range("R7").Value = 0.04741379310344827
$ws.Range("S7").Value = 0.3836206896551724
$ws.Range("B8").Value = 0.1121951219512195
$ws.Range("D8").Value = 0.01951219512195122
$ws.Range("F8").Value = 0.07073170731707316
$ws.Range("J8").Value = 0.1268292682926829
$ws.Range("O8").Value = 0.01951219512195122
$ws.Range("Q8").Value = 0.1731707317073171
$ws.Range("R8").Value = 0.08292682926829269
$ws.Range("S8").Value = 0.3951219512195122
$ws.Range("B9").Value = 0.08583690987124463
$ws.Range("D9").Value = 0.02145922746781116
$ws.Range("E9").Value = 0.004291845493562232
$ws.Range("F9").Value = 0.05579399141630902
$ws.Range("J9").Value = 0.111587982832618
$ws.Range("O9").Value = 0.01716738197424893
$ws.Range("Q9").Value = 0.1759656652360515
$ws.Range("R9").Value = 0.1030042918454936
$ws.Range("S9").Value = 0.4248927038626609
$ws.Range("B10").Value = 0.1229796205200281
$ws.Range("D10").Value = 0.01546029515108925
$ws.Range("E10").Value = 0.002108222066057625
$ws.Range("F10").Value = 0.08292340126493324
$ws.Range("J10").Value = 0.1314125087842586
$ws.Range("O10").Value = 0.02600140548137737
$ws.Range("Q10").Value = 0.185523541813071
$ws.Range("R10").Value = 0.07308503162333099
$ws.Range("S10").Value = 0.3605059732958538
$ws.Range("G11").Value = 0.1361111111111111
$ws.Range("J11").Value = 0.075
$ws.Range("K11").Value = 0.1777777777777778
$ws.Range("L11").Value = 0.6027777777777777
$ws.Range("S11").Value = 0.008333333333333333
$ws.Range("G12").Value = 0.7521739130434782
$ws.Range("J12").Value = 0.1782608695652174
$ws.Range("K12").Value = 0.02608695652173913
$ws.Range("L12").Value = 0.02173913043478261
$ws.Range("S12").Value = 0.02173913043478261
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.375
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.05109489051094891
$ws.Range("H15").Value = 0.1496350364963504
$ws.Range("I15").Value = 0.0583941605839416
$ws.Range("J15").Value = 0.3394160583941606
$ws.Range("K15").Value = 0.0583941605839416
$ws.Range("M15").Value = 0.0072992700729927
$ws.Range("O15").Value = 0.06204379562043796
$ws.Range("S15").Value = 0.2737226277372263
$ws.Range("F16").Value = 0.00816326530612245
$ws.Range("H16").Value = 0.1265306122448979
$ws.Range("I16").Value = 0.07346938775510205
$ws.Range("J16").Value = 0.436734693877551
$ws.Range("K16").Value = 0.1469387755102041
$ws.Range("M16").Value = 0.0163265306122449
$ws.Range("O16").Value = 0.06122448979591837
$ws.Range("S16").Value = 0.1306122448979592
$ws.Range("F17").Value = 0.02838427947598253
$ws.Range("H17").Value = 0.1310043668122271
$ws.Range("I17").Value = 0.09388646288209607
$ws.Range("J17").Value = 0.388646288209607
$ws.Range("K17").Value = 0.1179039301310044
$ws.Range("M17").Value = 0.01746724890829694
$ws.Range("O17").Value = 0.08078602620087336
$ws.Range("S17").Value = 0.1419213973799127
$ws.Range("F18").Value = 0.02870813397129187
$ws.Range("H18").Value = 0.1339712918660287
$ws.Range("I18").Value = 0.08133971291866028
$ws.Range("J18").Value = 0.4114832535885167
$ws.Range("K18").Value = 0.08133971291866028
$ws.Range("M18").Value = 0.01913875598086124
$ws.Range("O18").Value = 0.1100478468899522
$ws.Range("S18").Value = 0.1339712918660287
$ws.Range("F19").Value = 0.02068473609129814
$ws.Range("H19").Value = 0.1761768901569187
$ws.Range("I19").Value = 0.1005706134094151
$ws.Range("J19").Value = 0.3666191155492154
$ws.Range("K19").Value = 0.1212553495007133
$ws.Range("M19").Value = 0.01997146932952924
$ws.Range("N19").Value = 0.0007132667617689016
$ws.Range("O19").Value = 0.06776034236804565
$ws.Range("S19").Value = 0.1262482168330956
